$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 132.91667
$ws.Range("I33").Value = 115.14286
$ws.Range("K33").Value = 115.14286
$ws.Range("M33").Value = 113.85714

$ws.Range("H38").Value = 153.6
$ws.Range("I38").Value = 153.6
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 460.8
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = -88.79999999999995
$ws.Range("M38").ClearContents()

$ws.Range("H40").Value = 3450
$ws.Range("I40").Value = 2300
$ws.Range("J40").Value = 3833.3333
$ws.Range("K40").Value = 2300
$ws.Range("L40").Value = 3833.3333
$ws.Range("M40").Value = -2125
$ws.Range("N40").Value = -4183.3333

$ws.Range("H43").Value = 4499.5
$ws.Range("I43").Value = 4499.5
$ws.Range("K43").Value = 4499.5
$ws.Range("M43").Value = -4430.5

$ws.Range("H58").Value = 2365.2144

$ws.Range("H98").Value = 1247.5555
$ws.Range("I98").Value = 1027.75
$ws.Range("K98").Value = 1027.75
$ws.Range("M98").Value = 470.25

$ws.Range("H122").Value = 1247.5555
$ws.Range("I122").Value = 1027.75
$ws.Range("K122").Value = 3083.25
$ws.Range("M122").Value = -633.25

$ws.Range("H138").Value = 2463.9412
$ws.Range("I138").Value = 1480.7273
$ws.Range("K138").Value = 4442.1819
$ws.Range("M138").Value = 697.8181000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 36375.5
$ws.Range("I31").Value = 10471
$ws.Range("J31").Value = 62280
$ws.Range("K31").Value = 10471
$ws.Range("L31").Value = 62280
$ws.Range("M31").Value = -10177
$ws.Range("N31").Value = -62868

$ws.Range("H32").Value = 4022.2144
$ws.Range("I32").Value = 3754.6924
$ws.Range("K32").Value = 3754.6924
$ws.Range("M32").Value = -3467.6924

$ws.Range("H97").Value = 898.0833
$ws.Range("I97").Value = 777.7
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 777.7
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -281.7
$ws.Range("N97").Value = -2492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7749
$ws.Range("I86").Value = 10000
$ws.Range("J86").Value = 7298.8
$ws.Range("K86").Value = 10000
$ws.Range("L86").Value = 7298.8
$ws.Range("M86").Value = -8877
$ws.Range("N86").Value = -9544.799999999999

$ws.Range("H89").Value = 7749
$ws.Range("I89").Value = 10000
$ws.Range("J89").Value = 7298.8
$ws.Range("K89").Value = 50000
$ws.Range("L89").Value = 36494
$ws.Range("M89").Value = -44384
$ws.Range("N89").Value = -47726

$ws.Range("H94").Value = 1603.3334
$ws.Range("I94").Value = 1508.3334
$ws.Range("J94").Value = 2173.3333
$ws.Range("K94").Value = 1508.3334
$ws.Range("L94").Value = 2173.3333
$ws.Range("M94").Value = -1057.3334
$ws.Range("N94").Value = -3075.3333

$ws.Range("H107").Value = 967
$ws.Range("I107").Value = 819.5714
$ws.Range("K107").Value = 819.5714
$ws.Range("M107").Value = 1100.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3447.7856
$ws.Range("I31").Value = 2494.5
$ws.Range("J31").Value = 4162.75
$ws.Range("K31").Value = 2494.5
$ws.Range("L31").Value = 4162.75
$ws.Range("M31").Value = -2199.5
$ws.Range("N31").Value = -4752.75

$ws.Range("H34").Value = 3447.7856
$ws.Range("I34").Value = 2494.5
$ws.Range("J34").Value = 4162.75
$ws.Range("K34").Value = 2494.5
$ws.Range("L34").Value = 4162.75
$ws.Range("M34").Value = -2292.5
$ws.Range("N34").Value = -4566.75

$ws.Range("H81").Value = 52069.25
$ws.Range("I81").Value = 15000
$ws.Range("J81").Value = 57364.855
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 57364.855
$ws.Range("N81").Value = -59360.855
$ws.Range("M81").Value = -14002

$ws.Range("H84").Value = 52069.25
$ws.Range("I84").Value = 15000
$ws.Range("J84").Value = 57364.855
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 172094.565
$ws.Range("N84").Value = -182078.565
$ws.Range("M84").Value = -40008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 289
$ws.Range("J23").Value = 356.25
$ws.Range("L23").Value = 1068.75
$ws.Range("N23").Value = -1538.75

$ws.Range("H97").Value = 451
$ws.Range("I97").Value = 250
$ws.Range("K97").Value = 750
$ws.Range("M97").Value = -254

$ws.Range("H111").Value = 93.666664
$ws.Range("I111").Value = 93.666664
$ws.Range("K111").Value = 280.999992
$ws.Range("M111").Value = 2786.000008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2987
$ws.Range("I102").Value = 2987
$ws.Range("K102").Value = 2987
$ws.Range("M102").Value = -1365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1340
$ws.Range("J22").Value = 1340
$ws.Range("L22").Value = 1340
$ws.Range("N22").Value = -1930

$ws.Range("H27").Value = 1340
$ws.Range("J27").Value = 1340
$ws.Range("L27").Value = 1340
$ws.Range("N27").Value = -1554

$ws.Range("H46").Value = 3916.5833
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 4299.9
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 4299.9
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -4675.9

$ws.Range("H94").Value = 60000
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31040

$ws.Range("H63").Value = 32249
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 32249
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 32249
$ws.Range("N63").Value = -33497
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 32249
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 32249
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 96747
$ws.Range("N66").Value = -102987
$ws.Range("M66").ClearContents()

$ws.Range("H81").Value = 2543.3333
$ws.Range("I81").Value = 2543.3333
$ws.Range("K81").Value = 5086.6666
$ws.Range("M81").Value = -4025.6666

$ws.Range("H84").Value = 2543.3333
$ws.Range("I84").Value = 2543.3333
$ws.Range("K84").Value = 25433.333
$ws.Range("M84").Value = -20129.333

$ws.Range("H122").Value = 4268
$ws.Range("I122").Value = 4268
$ws.Range("K122").Value = 12804
$ws.Range("M122").Value = -10354

$ws.Range("H126").Value = 1186.3334
$ws.Range("I126").Value = 779.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2338.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = 131.5
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 1294.0588
$ws.Range("I132").Value = 1280.2667
$ws.Range("K132").Value = 3840.800099999999
$ws.Range("M132").Value = -1310.800099999999

$ws.Range("H136").Value = 4309.355
$ws.Range("I136").Value = 4408.4546
$ws.Range("K136").Value = 13225.3638
$ws.Range("M136").Value = -10675.3638
